$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price/hour values are stored as literal text
# (matching the source data which stores everything as inline strings),
# by pre-formatting the D and G columns as Text before assigning values.
$ws.Range("D2:D51").NumberFormat = "@"
$ws.Range("G2:G51").NumberFormat = "@"

$ws.Range("D2").Value = '242.76'
$ws.Range("G2").Value = '20'
$ws.Range("D3").Value = '23.08'
$ws.Range("G3").Value = '20'
$ws.Range("D4").Value = '5.413'
$ws.Range("G4").Value = '20'
$ws.Range("D5").Value = '0.05890'
$ws.Range("G5").Value = '20'
$ws.Range("D6").Value = '3.430'
$ws.Range("G6").Value = '20'
$ws.Range("D7").Value = '6.532'
$ws.Range("G7").Value = '20'
$ws.Range("D8").Value = '0.8085'
$ws.Range("G8").Value = '20'
$ws.Range("D9").Value = '0.9310'
$ws.Range("G9").Value = '20'
$ws.Range("D10").Value = '0.1418'
$ws.Range("G10").Value = '20'
$ws.Range("G11").Value = '20'
$ws.Range("D12").Value = '0.03312'
$ws.Range("G12").Value = '20'
$ws.Range("B13").Value = 'ProBitToken'
$ws.Range("C13").Value = 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'
$ws.Range("D13").Value = '0.1331'
$ws.Range("E13").Value = '12ProBitTokenPROB'
$ws.Range("G13").Value = '20'
$ws.Range("B14").Value = 'BitrueCoin'
$ws.Range("C14").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D14").Value = '0.03058'
$ws.Range("E14").Value = '13BitrueCoinBTR'
$ws.Range("G14").Value = '20'
$ws.Range("B15").Value = 'BitMartToken'
$ws.Range("C15").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D15").Value = '0.09346'
$ws.Range("E15").Value = '14BitMartTokenBMX'
$ws.Range("G15").Value = '20'
$ws.Range("B16").Value = 'MCDex'
$ws.Range("C16").Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range("D16").Value = '3.855'
$ws.Range("E16").Value = '15MCDexMCB'
$ws.Range("G16").Value = '20'
$ws.Range("B17").Value = 'BitForexToken'
$ws.Range("C17").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D17").Value = '0.001575'
$ws.Range("E17").Value = '16BitForexTokenBF'
$ws.Range("G17").Value = '20'
$ws.Range("B18").Value = 'CoinExToken'
$ws.Range("C18").Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range("D18").Value = '0.04671'
$ws.Range("E18").Value = '17CoinExTokenCET'
$ws.Range("G18").Value = '20'
$ws.Range("B19").Value = 'One'
$ws.Range("C19").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D19").Value = '0.0005900'
$ws.Range("E19").Value = '18OneONE'
$ws.Range("G19").Value = '20'
$ws.Range("B20").Value = 'TigerCash'
$ws.Range("C20").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D20").Value = '0.005881'
$ws.Range("E20").Value = '19TigerCashTCH'
$ws.Range("G20").Value = '20'
$ws.Range("B21").Value = 'BitKan'
$ws.Range("C21").Value = 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'
$ws.Range("D21").Value = '0.001264'
$ws.Range("E21").Value = '20BitKanKAN'
$ws.Range("G21").Value = '20'
$ws.Range("B22").Value = 'HotbitToken'
$ws.Range("C22").Value = 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'
$ws.Range("D22").Value = '0.004901'
$ws.Range("E22").Value = '21HotbitTokenHTB'
$ws.Range("G22").Value = '20'
$ws.Range("B23").Value = 'NitroEx'
$ws.Range("C23").Value = 'https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx'
$ws.Range("D23").Value = '0.00006802'
$ws.Range("E23").Value = '22NitroExNTX'
$ws.Range("G23").Value = '20'
$ws.Range("B24").Value = 'LEO'
$ws.Range("C24").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D24").Value = '3.564'
$ws.Range("E24").Value = '23LEOLEO'
$ws.Range("G24").Value = '20'
$ws.Range("B25").Value = 'BTSEToken'
$ws.Range("C25").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D25").Value = '2.147'
$ws.Range("E25").Value = '24BTSETokenBTSE'
$ws.Range("G25").Value = '20'
$ws.Range("B26").Value = 'BitpandaEcosystemToken'
$ws.Range("C26").Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range("D26").Value = '0.3233'
$ws.Range("E26").Value = '25BitpandaEcosystemTokenBEST'
$ws.Range("G26").Value = '20'
$ws.Range("D27").Value = '0.0002296'
$ws.Range("G27").Value = '20'
$ws.Range("G28").Value = '20'
$ws.Range("G29").Value = '20'
$ws.Range("G30").Value = '20'
$ws.Range("G31").Value = '20'
$ws.Range("G32").Value = '20'
$ws.Range("G33").Value = '20'
$ws.Range("G34").Value = '20'
$ws.Range("G35").Value = '20'
$ws.Range("G36").Value = '20'
$ws.Range("G37").Value = '20'
$ws.Range("G38").Value = '20'
$ws.Range("G39").Value = '20'
$ws.Range("D40").Value = '0.03974'
$ws.Range("G40").Value = '20'
$ws.Range("D41").Value = '0.006180'
$ws.Range("G41").Value = '20'
$ws.Range("D42").Value = '0.1074'
$ws.Range("G42").Value = '20'
$ws.Range("D43").Value = '0.002571'
$ws.Range("G43").Value = '20'
$ws.Range("D44").Value = '0.009683'
$ws.Range("G44").Value = '20'
$ws.Range("D45").Value = '0.00005176'
$ws.Range("G45").Value = '20'
$ws.Range("G46").Value = '20'
$ws.Range("D47").Value = '0.6700'
$ws.Range("G47").Value = '20'
$ws.Range("G48").Value = '20'
$ws.Range("D49").Value = '0.00002100'
$ws.Range("G49").Value = '20'
$ws.Range("D50").Value = '0.0002000'
$ws.Range("G50").Value = '20'
$ws.Range("G51").Value = '20'
